$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.841.27'
$ws.Cells.Item(2, 5).Value = '  +1.00%  '
$ws.Cells.Item(3, 4).Value = '2.088.13'
$ws.Cells.Item(3, 5).Value = '  +0.68%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '235.05'
$ws.Cells.Item(5, 5).Value = '  -0.17%  '
$ws.Cells.Item(6, 5).Value = '  -0.29%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '59.50'
$ws.Cells.Item(7, 5).Value = '  +3.55%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '1.00'
$ws.Cells.Item(8, 5).Value = '  -0.09%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.392'
$ws.Cells.Item(9, 5).Value = '  -0.82%  '
$ws.Cells.Item(10, 5).Value = '  +2.37%  '
$ws.Cells.Item(11, 5).Value = '  +2.86%  '
$ws.Cells.Item(12, 4).Value = '2.394.38'
$ws.Cells.Item(12, 5).Value = '  +0.68%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '14.74'
$ws.Cells.Item(13, 5).Value = '  +1.55%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '21.42'
$ws.Cells.Item(14, 5).Value = '  +3.58%  '
$ws.Cells.Item(15, 5).Value = '  -0.97%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '5.34'
$ws.Cells.Item(16, 5).Value = '  +2.75%  '
$ws.Cells.Item(17, 4).Value = '2.108.98'
$ws.Cells.Item(17, 5).Value = '  +1.73%  '
$ws.Cells.Item(18, 4).Value = '37.779.34'
$ws.Cells.Item(18, 5).Value = '  +1.03%  '
$ws.Cells.Item(19, 5).Value = '  -3.53%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '71.66'
$ws.Cells.Item(21, 4).Value = '0.0₃0830'
$ws.Cells.Item(21, 5).Value = '  +1.20%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '229.02'
$ws.Cells.Item(22, 5).Value = '  +0.81%  '
$ws.Cells.Item(23, 5).Value = '  -0.03%  '
$ws.Cells.Item(24, 5).Value = '  -0.66%  '
$ws.Cells.Item(25, 5).Value = '  +0.66%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '170.52'
$ws.Cells.Item(26, 5).Value = '  +2.00%  '
$ws.Cells.Item(27, 5).Value = '  +8.08%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.07'
$ws.Cells.Item(28, 5).Value = '  +2.35%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.44'
$ws.Cells.Item(29, 5).Value = '  -0.28%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '19.61'
$ws.Cells.Item(30, 5).Value = '  +2.32%  '
$ws.Cells.Item(31, 5).Value = '  +1.79%  '
$ws.Cells.Item(32, 5).Value = '  +3.54%  '
$ws.Cells.Item(33, 5).Value = '  +1.77%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.71'
$ws.Cells.Item(34, 5).Value = '  +2.91%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.54'
$ws.Cells.Item(35, 5).Value = '  +0.97%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '3.54'
$ws.Cells.Item(36, 5).Value = '  +6.46%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.999'
$ws.Cells.Item(38, 5).Value = '  -0.11%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.46'
$ws.Cells.Item(39, 5).Value = '  -4.33%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0987'
$ws.Cells.Item(40, 5).Value = '  +2.01%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '100.08'
$ws.Cells.Item(41, 5).Value = '  +1.34%  '
$ws.Cells.Item(42, 5).Value = '  +0.03%  '
$ws.Cells.Item(43, 5).Value = '  +0.62%  '
$ws.Cells.Item(44, 4).Value = '1.463.99'
$ws.Cells.Item(44, 5).Value = '  -1.21%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.17'
$ws.Cells.Item(45, 5).Value = '  -0.56%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '4.21'
$ws.Cells.Item(46, 5).Value = '  +2.45%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '16.21'
$ws.Cells.Item(47, 5).Value = '  +5.48%  '
$ws.Cells.Item(48, 5).Value = '  +3.71%  '
$ws.Cells.Item(49, 5).Value = '  +2.83%  '
$ws.Cells.Item(50, 5).Value = '  +2.52%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '47.69'
$ws.Cells.Item(51, 5).Value = '  +6.41%  '
